$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2826.6365
$ws.Range("J32").Value = 2677
$ws.Range("L32").Value = 2677
$ws.Range("N32").Value = -3329

$ws.Range("H62").Value = 2969.348
$ws.Range("I62").Value = 2914.3333
$ws.Range("J62").Value = 3029.3635
$ws.Range("K62").Value = 2914.3333
$ws.Range("L62").Value = 3029.3635
$ws.Range("M62").Value = -2290.3333
$ws.Range("N62").Value = -4277.363499999999

$ws.Range("H65").Value = 2969.348
$ws.Range("I65").Value = 2914.3333
$ws.Range("J65").Value = 3029.3635
$ws.Range("K65").Value = 14571.6665
$ws.Range("L65").Value = 15146.8175
$ws.Range("M65").Value = -11451.6665
$ws.Range("N65").Value = -21386.8175

$ws.Range("H98").Value = 890.7143
$ws.Range("I98").Value = 552.2222
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 552.2222
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 945.7778
$ws.Range("N98").Value = -4496

$ws.Range("H113").Value = 1538
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1538
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1538
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8046

$ws.Range("H116").Value = 2082242.2
$ws.Range("I116").Value = 7694637
$ws.Range("J116").Value = 3577.5557
$ws.Range("K116").Value = 7694637
$ws.Range("L116").Value = 3577.5557
$ws.Range("M116").Value = -7691195
$ws.Range("N116").Value = -10461.5557

$ws.Range("H122").Value = 890.7143
$ws.Range("I122").Value = 552.2222
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 1656.6666
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = 793.3334
$ws.Range("N122").Value = -9400

$ws.Range("H137").Value = 2139.2144
$ws.Range("I137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("M137").Value = -450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1500.909
$ws.Range("I2").Value = 1478.1765
$ws.Range("J2").Value = 1578.2
$ws.Range("K2").Value = 1478.1765
$ws.Range("L2").Value = 1578.2
$ws.Range("M2").Value = -1365.1765
$ws.Range("N2").Value = -1804.2

$ws.Range("H61").Value = 2590.5
$ws.Range("I61").Value = 1704.8
$ws.Range("J61").Value = 4066.6667
$ws.Range("K61").Value = 1704.8
$ws.Range("L61").Value = 4066.6667
$ws.Range("M61").Value = -1492.8
$ws.Range("N61").Value = -4490.6667

$ws.Range("H74").Value = 1520.2413
$ws.Range("I74").Value = 1608.3684
$ws.Range("K74").Value = 1608.3684
$ws.Range("M74").Value = -734.3684000000001

$ws.Range("H77").Value = 1520.2413
$ws.Range("I77").Value = 1608.3684
$ws.Range("K77").Value = 8041.842000000001
$ws.Range("M77").Value = -3673.842000000001

$ws.Range("H92").Value = 29966.666
$ws.Range("J92").Value = 29966.666
$ws.Range("L92").Value = 29966.666
$ws.Range("N92").Value = -34958.666

$ws.Range("H116").Value = 1500.909
$ws.Range("I116").Value = 1478.1765
$ws.Range("J116").Value = 1578.2
$ws.Range("K116").Value = 1478.1765
$ws.Range("L116").Value = 1578.2
$ws.Range("M116").Value = 815.8235
$ws.Range("N116").Value = -6166.2

$ws.Range("H136").Value = 2590.5
$ws.Range("I136").Value = 1704.8
$ws.Range("J136").Value = 4066.6667
$ws.Range("K136").Value = 5114.4
$ws.Range("L136").Value = 12200.0001
$ws.Range("M136").Value = -2564.4
$ws.Range("N136").Value = -17300.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1500.909
$ws.Range("I3").Value = 1478.1765
$ws.Range("J3").Value = 1578.2
$ws.Range("K3").Value = 1478.1765
$ws.Range("L3").Value = 1578.2
$ws.Range("M3").Value = -1364.1765
$ws.Range("N3").Value = -1806.2

$ws.Range("H134").Value = 2294.3462
$ws.Range("I134").Value = 2192
$ws.Range("J134").Value = 2635.5
$ws.Range("K134").Value = 6576
$ws.Range("L134").Value = 7906.5
$ws.Range("M134").Value = -4041
$ws.Range("N134").Value = -12976.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3800
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3800
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3800
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -4024

$ws.Range("H132").Value = 3159.4348
$ws.Range("I132").Value = 1843.7333
$ws.Range("J132").Value = 5626.375
$ws.Range("K132").Value = 5531.199900000001
$ws.Range("L132").Value = 16879.125
$ws.Range("M132").Value = -3001.199900000001
$ws.Range("N132").Value = -21939.125

$ws.Range("H134").Value = 3910.6667
$ws.Range("I134").Value = 2057.4546
$ws.Range("J134").Value = 9007
$ws.Range("K134").Value = 6172.3638
$ws.Range("L134").Value = 27021
$ws.Range("M134").Value = -3637.3638
$ws.Range("N134").Value = -32091

$ws.Range("H135").Value = 24913
$ws.Range("J135").Value = 24913
$ws.Range("L135").Value = 24913
$ws.Range("N135").Value = -35053

$ws.Range("H138").Value = 39343.25
$ws.Range("J138").Value = 39343.25
$ws.Range("L138").Value = 39343.25
$ws.Range("N138").Value = -49623.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3466.5
$ws.Range("I94").Value = 899.5
$ws.Range("J94").Value = 4750
$ws.Range("K94").Value = 2698.5
$ws.Range("L94").Value = 14250
$ws.Range("M94").Value = -2022.5
$ws.Range("N94").Value = -15602

$ws.Range("H98").Value = 718.4167
$ws.Range("J98").Value = 809
$ws.Range("L98").Value = 2427
$ws.Range("N98").Value = -5423

$ws.Range("H99").Value = 1728.625
$ws.Range("I99").Value = 1204.4615
$ws.Range("K99").Value = 3613.3845
$ws.Range("M99").Value = -1367.3845

$ws.Range("H101").Value = 4500
$ws.Range("I101").Value = 3000
$ws.Range("K101").Value = 9000
$ws.Range("M101").Value = -6566

$ws.Range("H107").Value = 242.72728
$ws.Range("I107").Value = 174.28572
$ws.Range("J107").Value = 362.5
$ws.Range("K107").Value = 522.85716
$ws.Range("L107").Value = 1087.5
$ws.Range("M107").Value = 1397.14284
$ws.Range("N107").Value = -4927.5

$ws.Range("H113").Value = 1724759.2
$ws.Range("I113").Value = 4310900.5
$ws.Range("J113").Value = 665
$ws.Range("K113").Value = 12932701.5
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = -12930531.5
$ws.Range("N113").Value = -6335

$ws.Range("H132").Value = 1780.8
$ws.Range("I132").Value = 904
$ws.Range("K132").Value = 8136
$ws.Range("M132").Value = -5606

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3707988
$ws.Range("I126").Value = 6670749.5
$ws.Range("J126").Value = 4536.3335
$ws.Range("K126").Value = 20012248.5
$ws.Range("L126").Value = 13609.0005
$ws.Range("M126").Value = -20009778.5
$ws.Range("N126").Value = -18549.0005

$ws.Range("H132").Value = 3253.3333
$ws.Range("I132").Value = 2833.5
$ws.Range("K132").Value = 8500.5
$ws.Range("M132").Value = -5970.5

$ws.Range("H141").Value = 47981
$ws.Range("J141").Value = 47981
$ws.Range("L141").Value = 47981
$ws.Range("N141").Value = -58341

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 54619.332
$ws.Range("J128").Value = 54619.332
$ws.Range("L128").Value = 54619.332
$ws.Range("N128").Value = -64579.332

$ws.Range("H136").Value = 3276.0833
$ws.Range("I136").Value = 2731.3
$ws.Range("K136").Value = 8193.900000000001
$ws.Range("M136").Value = -5643.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 420
$ws.Range("I107").Value = 382.85715
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 1148.57145
$ws.Range("L107").Value = 1650
$ws.Range("M107").Value = 771.4285500000001
$ws.Range("N107").Value = -5490
